# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (new fund-holdings data) right after the
# "总计" (summary) sheet and before the existing "2021-Q3" sheet, then
# updates the "总计" sheet so that row 2 now reports the 2022-Q3 totals and a
# new row 3 carries over the original 2021-Q3 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (i.e. before the
#    existing "2021-Q3" sheet), mirroring the workbook.xml sheet order.
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$q3_2022 = $wb.Worksheets.Add($null, $summarySheet)
$q3_2022.Name = "2022-Q3"

# Header row (B1:H1) - text labels matching the "2021-Q3" sheet layout.
$q3_2022.Range("B1").Value = "基金代码"
$q3_2022.Range("C1").Value = "基金名称"
$q3_2022.Range("D1").Value = "基金规模"
$q3_2022.Range("E1").Value = "股票总仓位"
$q3_2022.Range("F1").Value = "仓位占比"
$q3_2022.Range("G1").Value = "持有市值(亿元)"
$q3_2022.Range("H1").Value = "仓位排名"

# Header cells use the same bold/bordered style as the "总计" sheet header
# (style index 2 in the original workbook).
$summarySheet.Range("B1").Copy()
$q3_2022.Range("B1:H1").PasteSpecial(-4122)

# Row 2 data
$q3_2022.Range("A2").Value = 0
$q3_2022.Range("B2").Value = "'016950"
$q3_2022.Range("B2").Style = "Normal"
$q3_2022.Range("C2").Value = "鹏华睿投灵活配置混合C"
$q3_2022.Range("D2").Value = "'4.12"
$q3_2022.Range("D2").Style = "Normal"
$q3_2022.Range("E2").Value = "'83.97"
$q3_2022.Range("E2").Style = "Normal"
$q3_2022.Range("F2").Value = "'1.83"
$q3_2022.Range("F2").Style = "Normal"
$q3_2022.Range("G2").Value = "'0.0754"
$q3_2022.Range("G2").Style = "Normal"
$q3_2022.Range("H2").Value = 9

# Row 3 data
$q3_2022.Range("A3").Value = 1
$q3_2022.Range("B3").Value = "'005434"
$q3_2022.Range("B3").Style = "Normal"
$q3_2022.Range("C3").Value = "鹏华睿投灵活配置混合A"
$q3_2022.Range("D3").Value = "'0.00"
$q3_2022.Range("D3").Style = "Normal"
$q3_2022.Range("E3").Value = "'83.97"
$q3_2022.Range("E3").Style = "Normal"
$q3_2022.Range("F3").Value = "'1.83"
$q3_2022.Range("F3").Style = "Normal"
$q3_2022.Range("G3").Value = 0
$q3_2022.Range("H3").Value = 9

# A2/A3 share the "总计"-header style too (matches s="2" in the target).
$summarySheet.Range("B1").Copy()
$q3_2022.Range("A2:A3").PasteSpecial(-4122)
$q3_2022.Range("A2").Value = 0
$q3_2022.Range("A3").Value = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: row 2 now reflects 2022-Q3, and a new row 3
#    carries the original 2021-Q3 totals forward.
# ---------------------------------------------------------------------
$summarySheet.Range("A2").Copy()
$summarySheet.Range("A3").PasteSpecial(-4122)

$summarySheet.Range("B2").Value = "2022-Q3"
$summarySheet.Range("D2").Value = 0.08

$summarySheet.Range("A3").Value = 1
$summarySheet.Range("B3").Value = "2021-Q3"
$summarySheet.Range("C3").Value = 2
$summarySheet.Range("D3").Value = 0.1
